$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column H header + data (Password column)
$ws.Range("H1").Value = "Password"
$ws.Range("H2:H24").Value = "123abc"

# Set column H width
$ws.Columns.Item(8).ColumnWidth = 13.28515625

# Style H1 like the rest of the header row (bold, centered, no border) -
# it already inherits the row-level formatting, but set alignment explicitly.
$ws.Range("H1").HorizontalAlignment = -4108  # xlCenter

# Style H2:H24 - add thin left/right border and center alignment
$dataRange = $ws.Range("H2:H24")
$dataRange.HorizontalAlignment = -4108  # xlCenter
$dataRange.Borders.Item(7).LineStyle = 1   # xlEdgeLeft = 7, xlContinuous = 1
$dataRange.Borders.Item(7).Weight = 2      # xlThin = 2
$dataRange.Borders.Item(10).LineStyle = 1  # xlEdgeRight = 10
$dataRange.Borders.Item(10).Weight = 2     # xlThin = 2

# Update selection to mimic final cursor position
$ws.Range("H32").Select()
